# Updated symbol list on Sat Dec 24 20:56:56 UTC 2022 with GitHub Actions
#
# Applies the refreshed coin price / ranking snapshot to Sheet1.
# Column D ("Price") values are written with a leading apostrophe so
# Excel stores them as literal text (matching the original inlineStr
# cells) instead of re-parsing them as numbers and dropping
# significant trailing zeros (e.g. "0.03400" -> 0.034).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- simple price refreshes (rows whose coin identity did not change) ---
$ws.Range("D2").Value  = "'244.56"
$ws.Range("D3").Value  = "'21.92"
$ws.Range("D4").Value  = "'5.402"
$ws.Range("D5").Value  = "'0.06034"
$ws.Range("D6").Value  = "'3.393"
$ws.Range("D7").Value  = "'0.8142"
$ws.Range("D8").Value  = "'0.9273"

# --- rows 9-17: coin ranking reshuffled, so Coin/Link/Price/Volume all change ---
$ws.Range("B9").Value  = "WazirX"
$ws.Range("C9").Value  = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value  = "'0.1435"
$ws.Range("E9").Value  = "8WazirXWRX"

$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").Value = "'0.07509"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.03400"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.03039"
$ws.Range("E12").Value = "11BitrueCoinBTR"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09418"
$ws.Range("E13").Value = "12BitMartTokenBMX"

$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D14").Value = "'4.011"
$ws.Range("E14").Value = "13MCDexMCB"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001589"
$ws.Range("E15").Value = "14BitForexTokenBF"

$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D16").Value = "'0.04822"
$ws.Range("E16").Value = "15CoinExTokenCET"

$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D17").Value = "'0.0005941"
$ws.Range("E17").Value = "16OneONE"

# --- more simple price refreshes ---
$ws.Range("D18").Value = "'0.005520"
$ws.Range("D19").Value = "'0.004165"
$ws.Range("D20").Value = "'0.0009874"
$ws.Range("D21").Value = "'3.667"
$ws.Range("D22").Value = "'6.427"
$ws.Range("D26").Value = "'0.00008402"
$ws.Range("D27").Value = "'0.0002900"
$ws.Range("D40").Value = "'0.03994"

# --- rows 41-43: coin ranking reshuffled again ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").Value = "'0.1077"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.002721"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").Value = "'0.003053"
$ws.Range("E43").Value = "42KickTokenKICK"

# --- trailing simple price refreshes + one "Best/Worst in 24h" label tweak ---
$ws.Range("D44").Value = "'0.005796"
$ws.Range("D45").Value = "'0.00005240"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOINBestin24h"
$ws.Range("D48").Value = "'0.002324"
$ws.Range("D49").Value = "'0.00002100"
